# Generate Report for Handback
# Adds a "version mismatch" handback-report row to both locale sheets
# (zh-cn and de-de): the Latest Target File / Latest Handback DateTime /
# Error Detail columns get populated and a hyperlink is added on the
# "Latest Target File" cell (column I), mirroring the existing hyperlink
# already present on column A for the same source doc.

$wb = $excel.ActiveWorkbook

$docName = "dad6e231-0e69-4f66-b3a3-a4e981dddb41.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/83d206115fd0722261fcb1544ae522720a505d23/e2e/dad6e231-0e69-4f66-b3a3-a4e981dddb41.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a74a1dffb33741bd41b7024f660af3513ef23333/e2e/dad6e231-0e69-4f66-b3a3-a4e981dddb41.md."
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a74a1dffb33741bd41b7024f660af3513ef23333/e2e/dad6e231-0e69-4f66-b3a3-a4e981dddb41.md"

# --- zh-cn sheet -----------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $docName)
$wsZh.Range("I8").Style = "HyperLink"

$wsZh.Range("J8").Value = "dad6e231-0e69-4f66-b3a3-a4e981dddb41.3b2c26de9ebd76b4deda9d93d7b232bccbd81edd.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-28 00:41:02"
$wsZh.Range("P8").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet -------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $docName)
$wsDe.Range("I8").Style = "HyperLink"

$wsDe.Range("J8").Value = "dad6e231-0e69-4f66-b3a3-a4e981dddb41.3b2c26de9ebd76b4deda9d93d7b232bccbd81edd.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-28 00:41:11"
$wsDe.Range("P8").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
